# Swap the stat values between row 2 and row 4 (columns C:F) on the
# active sheet. Values are stored as text (not numbers), so the
# number format is forced to "@" (Text) before assignment to keep
# the cells' underlying type as text. Reads use .Value2 (plain read
# of the stored value) since .Value is used here only for writes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("C", "D", "E", "F")

# Capture current values for row 2 and row 4 before overwriting anything.
$row2 = @{}
$row4 = @{}
foreach ($col in $cols) {
    $addr2 = $col + "2"
    $addr4 = $col + "4"
    $row2[$col] = $ws.Range($addr2).Value2
    $row4[$col] = $ws.Range($addr4).Value2
}

# Write row 2 <- old row 4, row 4 <- old row 2, preserving text storage.
foreach ($col in $cols) {
    $addr2 = $col + "2"
    $addr4 = $col + "4"

    $ws.Range($addr2).NumberFormat = "@"
    $ws.Range($addr2).Value = $row4[$col]

    $ws.Range($addr4).NumberFormat = "@"
    $ws.Range($addr4).Value = $row2[$col]
}
